# Refactor currency conversion sheet: explicit source and target amounts.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("currency_conversions")

# --- Header row (row 1) ---------------------------------------------------
# Old layout: date | foreign_amount | source_fees | source_currency | target_currency | comment
# New layout: date | source_amount  | source_fees | source_currency | target_amount | target_fees | target_currency | comment

# Shift the old E1 (target_currency) and F1 (comment) out of the way first,
# so nothing gets clobbered while we rebuild the row left-to-right.
$ws.Range("H1").Value = $ws.Range("F1").Value()
$ws.Range("H1").Style = $ws.Range("F1").Style
$ws.Range("G1").Value = $ws.Range("E1").Value()
$ws.Range("G1").Style = $ws.Range("E1").Style

# New columns
$ws.Range("B1").Value = "source_amount"
$ws.Range("C1").Value = "source_fees"
$ws.Range("D1").Value = "source_currency"
$ws.Range("E1").Value = "target_amount"
$ws.Range("F1").Value = "target_fees"

# Match styling of the rest of the header row, then make the new first
# data column (source_amount) bold to set it apart.
$ws.Range("C1:G1").Style = $ws.Range("A1").Style
$ws.Range("B1").Style = $ws.Range("A1").Style
$ws.Range("B1").Font.Bold = $true

# --- Data row (row 2) ------------------------------------------------------
# Old layout: date | 150 (foreign_amount) | 0 (source_fees) | EUR | USD
# New layout: date | -1 (source_amount) | 0 (source_fees) | EUR | 150 (target_amount) | 0 (target_fees) | USD
$ws.Range("G2").Value = $ws.Range("E2").Value()
$ws.Range("D2").Value = $ws.Range("D2").Value()
$ws.Range("E2").Value = $ws.Range("B2").Value()
$ws.Range("B2").Value = -1
$ws.Range("F2").Value = 0

# Activate this sheet (was previously on buy_orders) -> sets workbookView activeTab
$ws.Activate()
